$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add a date value in cell C1, formatted as a short date (matches the
# new numFmtId="14" cellXfs entry introduced by the edit).
$cell = $ws.Range("C1")
$cell.Value = 44307
$cell.NumberFormat = "mm-dd-yy"
